# Apply updated TPM-derived values (NATMI LR-pairs) for Il16-Kcnj10 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "MuSCs"
$ws.Range("G2").Value = 3.186175
$ws.Range("H2").Value = 6.37235
$ws.Range("I2").Value = 0.1401839069896986
$ws.Range("J2").Value = 0.09997212951453034
$ws.Range("M2").Value = 0.3477795
$ws.Range("N2").Value = 0.695559
$ws.Range("O2").Value = 0.7300808741333977
$ws.Range("P2").Value = 0.6432659787865084
$ws.Range("Q2").Value = 1.1080863484125
$ws.Range("R2").Value = 4.43234539365
$ws.Range("S2").Value = 0.1023455893544741
$ws.Range("T2").Value = 0.06430866974353594
$ws.Range("D3").Value = "Neutrophils"
$ws.Range("G3").Value = 3.186175
$ws.Range("H3").Value = 6.37235
$ws.Range("I3").Value = 0.1401839069896986
$ws.Range("J3").Value = 0.09997212951453034
$ws.Range("M3").Value = 0.04758299999999999
$ws.Range("N3").Value = 0.142749
$ws.Range("O3").Value = 0.09988926384070786
$ws.Range("P3").Value = 0.1320169463780862
$ws.Range("Q3").Value = 0.151607765025
$ws.Range("R3").Value = 0.9096465901499999
$ws.Range("S3").Value = 0.01400286727151525
$ws.Range("T3").Value = 0.01319801526142284
$ws.Range("G4").Value = 3.186175
$ws.Range("H4").Value = 6.37235
$ws.Range("I4").Value = 0.1401839069896986
$ws.Range("J4").Value = 0.09997212951453034
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.080995
$ws.Range("N4").Value = 0.242985
$ws.Range("O4").Value = 0.1700298620258944
$ws.Range("P4").Value = 0.2247170748354054
$ws.Range("Q4").Value = 0.258064244125
$ws.Range("R4").Value = 1.54838546475
$ws.Range("S4").Value = 0.02383545036370927
$ws.Range("T4").Value = 0.02246544450957155
$ws.Range("D5").Value = "MuSCs"
$ws.Range("I5").Value = 0.0583987308700193
$ws.Range("J5").Value = 0.06247056753580275
$ws.Range("M5").Value = 0.3477795
$ws.Range("N5").Value = 0.695559
$ws.Range("O5").Value = 0.7300808741333977
$ws.Range("P5").Value = 0.6432659787865084
$ws.Range("Q5").Value = 0.4616138744545001
$ws.Range("R5").Value = 2.769683246727
$ws.Range("S5").Value = 0.04263579648186473
$ws.Range("T5").Value = 0.04018519077126683
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("I6").Value = 0.0583987308700193
$ws.Range("J6").Value = 0.06247056753580275
$ws.Range("M6").Value = 0.04758299999999999
$ws.Range("N6").Value = 0.142749
$ws.Range("O6").Value = 0.09988926384070786
$ws.Range("P6").Value = 0.1320169463780862
$ws.Range("Q6").Value = 0.063157756533
$ws.Range("R6").Value = 0.5684198087969999
$ws.Range("S6").Value = 0.005833406235837848
$ws.Range("T6").Value = 0.008247173564582686
$ws.Range("I7").Value = 0.0583987308700193
$ws.Range("J7").Value = 0.06247056753580275
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.080995
$ws.Range("N7").Value = 0.242985
$ws.Range("O7").Value = 0.1700298620258944
$ws.Range("P7").Value = 0.2247170748354054
$ws.Range("Q7").Value = 0.1075060944116667
$ws.Range("R7").Value = 0.9675548497050001
$ws.Range("S7").Value = 0.009929528152316723
$ws.Range("T7").Value = 0.01403820319995323
$ws.Range("D8").Value = "MuSCs"
$ws.Range("G8").Value = 3.568381333333333
$ws.Range("H8").Value = 10.705144
$ws.Range("I8").Value = 0.1570000508245079
$ws.Range("J8").Value = 0.1679468394610618
$ws.Range("M8").Value = 0.3477795
$ws.Range("N8").Value = 0.695559
$ws.Range("O8").Value = 0.7300808741333977
$ws.Range("P8").Value = 0.6432659787865084
$ws.Range("Q8").Value = 1.241009875916
$ws.Range("R8").Value = 7.446059255496001
$ws.Range("S8").Value = 0.1146227343449446
$ws.Range("T8").Value = 0.1080344880700205
$ws.Range("D9").Value = "Neutrophils"
$ws.Range("G9").Value = 3.568381333333333
$ws.Range("H9").Value = 10.705144
$ws.Range("I9").Value = 0.1570000508245079
$ws.Range("J9").Value = 0.1679468394610618
$ws.Range("M9").Value = 0.04758299999999999
$ws.Range("N9").Value = 0.142749
$ws.Range("O9").Value = 0.09988926384070786
$ws.Range("P9").Value = 0.1320169463780862
$ws.Range("Q9").Value = 0.169794288984
$ws.Range("R9").Value = 1.528148600856
$ws.Range("S9").Value = 0.01568261949981382
$ws.Range("T9").Value = 0.02217182889950006
$ws.Range("G10").Value = 3.568381333333333
$ws.Range("H10").Value = 10.705144
$ws.Range("I10").Value = 0.1570000508245079
$ws.Range("J10").Value = 0.1679468394610618
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.080995
$ws.Range("N10").Value = 0.242985
$ws.Range("O10").Value = 0.1700298620258944
$ws.Range("P10").Value = 0.2247170748354054
$ws.Range("Q10").Value = 0.2890210460933333
$ws.Range("R10").Value = 2.60118941484
$ws.Range("S10").Value = 0.0266946969797495
$ws.Range("T10").Value = 0.03774052249154124
$ws.Range("D11").Value = "MuSCs"
$ws.Range("G11").Value = 1.2581685
$ws.Range("H11").Value = 2.516337
$ws.Range("I11").Value = 0.0553563366674362
$ws.Range("J11").Value = 0.03947736211385199
$ws.Range("M11").Value = 0.3477795
$ws.Range("N11").Value = 0.695559
$ws.Range("O11").Value = 0.7300808741333977
$ws.Range("P11").Value = 0.6432659787865084
$ws.Range("Q11").Value = 0.43756521184575
$ws.Range("R11").Value = 1.750260847383
$ws.Range("S11").Value = 0.04041460266298448
$ws.Range("T11").Value = 0.02539444398007643
$ws.Range("D12").Value = "Neutrophils"
$ws.Range("G12").Value = 1.2581685
$ws.Range("H12").Value = 2.516337
$ws.Range("I12").Value = 0.0553563366674362
$ws.Range("J12").Value = 0.03947736211385199
$ws.Range("M12").Value = 0.04758299999999999
$ws.Range("N12").Value = 0.142749
$ws.Range("O12").Value = 0.09988926384070786
$ws.Range("P12").Value = 0.1320169463780862
$ws.Range("Q12").Value = 0.05986743173549999
$ws.Range("R12").Value = 0.359204590413
$ws.Range("S12").Value = 0.005529503718628585
$ws.Range("T12").Value = 0.005211680797332691
$ws.Range("G13").Value = 1.2581685
$ws.Range("H13").Value = 2.516337
$ws.Range("I13").Value = 0.0553563366674362
$ws.Range("J13").Value = 0.03947736211385199
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.080995
$ws.Range("N13").Value = 0.242985
$ws.Range("O13").Value = 0.1700298620258944
$ws.Range("P13").Value = 0.2247170748354054
$ws.Range("Q13").Value = 0.1019053576575
$ws.Range("R13").Value = 0.6114321459450001
$ws.Range("S13").Value = 0.009412230285823138
$ws.Range("T13").Value = 0.008871237336442874
$ws.Range("D14").Value = "MuSCs"
$ws.Range("G14").Value = 11.59014566666667
$ws.Range("H14").Value = 34.770437
$ws.Range("I14").Value = 0.5099380611965939
$ws.Range("J14").Value = 0.545493362894508
$ws.Range("M14").Value = 0.3477795
$ws.Range("N14").Value = 0.695559
$ws.Range("O14").Value = 0.7300808741333977
$ws.Range("P14").Value = 0.6432659787865084
$ws.Range("Q14").Value = 4.030815064880501
$ws.Range("R14").Value = 24.184890389283
$ws.Range("S14").Value = 0.3722960254722993
$ws.Range("T14").Value = 0.3508973220038797
$ws.Range("D15").Value = "Neutrophils"
$ws.Range("G15").Value = 11.59014566666667
$ws.Range("H15").Value = 34.770437
$ws.Range("I15").Value = 0.5099380611965939
$ws.Range("J15").Value = 0.545493362894508
$ws.Range("M15").Value = 0.04758299999999999
$ws.Range("N15").Value = 0.142749
$ws.Range("O15").Value = 0.09988926384070786
$ws.Range("P15").Value = 0.1320169463780862
$ws.Range("Q15").Value = 0.5514939012569999
$ws.Range("R15").Value = 4.963445111313
$ws.Range("S15").Value = 0.0509373375372856
$ws.Range("T15").Value = 0.0720143680388462
$ws.Range("G16").Value = 11.59014566666667
$ws.Range("H16").Value = 34.770437
$ws.Range("I16").Value = 0.5099380611965939
$ws.Range("J16").Value = 0.545493362894508
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.080995
$ws.Range("N16").Value = 0.242985
$ws.Range("O16").Value = 0.1700298620258944
$ws.Range("P16").Value = 0.2247170748354054
$ws.Range("Q16").Value = 0.9387438482716666
$ws.Range("R16").Value = 8.448694634445001
$ws.Range("S16").Value = 0.08670469818700897
$ws.Range("T16").Value = 0.1225816728517821
$ws.Range("D17").Value = "MuSCs"
$ws.Range("G17").Value = 1.798348
$ws.Range("H17").Value = 5.395044
$ws.Range("I17").Value = 0.07912291345174399
$ws.Range("J17").Value = 0.0846397384802451
$ws.Range("M17").Value = 0.3477795
$ws.Range("N17").Value = 0.695559
$ws.Range("O17").Value = 0.7300808741333977
$ws.Range("P17").Value = 0.6432659787865084
$ws.Range("Q17").Value = 0.625428568266
$ws.Range("R17").Value = 3.752571409596
$ws.Range("S17").Value = 0.05776612581683043
$ws.Range("T17").Value = 0.05444586421772896
$ws.Range("D18").Value = "Neutrophils"
$ws.Range("G18").Value = 1.798348
$ws.Range("H18").Value = 5.395044
$ws.Range("I18").Value = 0.07912291345174399
$ws.Range("J18").Value = 0.0846397384802451
$ws.Range("M18").Value = 0.04758299999999999
$ws.Range("N18").Value = 0.142749
$ws.Range("O18").Value = 0.09988926384070786
$ws.Range("P18").Value = 0.1320169463780862
$ws.Range("Q18").Value = 0.08557079288399999
$ws.Range("R18").Value = 0.770137135956
$ws.Range("S18").Value = 0.007903529577626748
$ws.Range("T18").Value = 0.01117387981640176
$ws.Range("G19").Value = 1.798348
$ws.Range("H19").Value = 5.395044
$ws.Range("I19").Value = 0.07912291345174399
$ws.Range("J19").Value = 0.0846397384802451
$ws.Range("K19").Value = 1
$ws.Range("L19").Value = 0.3333333333333333
$ws.Range("M19").Value = 0.080995
$ws.Range("N19").Value = 0.242985
$ws.Range("O19").Value = 0.1700298620258944
$ws.Range("P19").Value = 0.2247170748354054
$ws.Range("Q19").Value = 0.14565719626
$ws.Range("R19").Value = 1.31091476634
$ws.Range("S19").Value = 0.01345325805728682
$ws.Range("T19").Value = 0.01901999444611438

Write-Output "Applied 234 cell updates"
